# Apply cryptocurrency price/volume updates from the Wed Jul 10 06:20:47 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17/18 also swap which coin (Coin + Link) occupies each rank position.
$ws.Range("D2").Value = "'59.082.16"
$ws.Range("E2").Value = "'  +2.96%  "
$ws.Range("D3").Value = "'3.107.91"
$ws.Range("E3").Value = "'  +1.03%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'524.01"
$ws.Range("E5").Value = "'  +1.89%  "
$ws.Range("D6").Value = "'144.17"
$ws.Range("E6").Value = "'  +1.87%  "
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("D8").Value = "'0.439"
$ws.Range("E8").Value = "'  +0.86%  "
$ws.Range("D9").Value = "'7.40"
$ws.Range("E9").Value = "'  +1.91%  "
$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "'  +1.07%  "
$ws.Range("D11").Value = "'0.384"
$ws.Range("E11").Value = "'  +3.05%  "
$ws.Range("D12").Value = "'3.637.21"
$ws.Range("E12").Value = "'  +1.11%  "
$ws.Range("E13").Value = "'  +1.31%  "
$ws.Range("D14").Value = "'27.12"
$ws.Range("E14").Value = "'  +5.93%  "
$ws.Range("E15").Value = "'  +1.44%  "
$ws.Range("D16").Value = "'59.021.21"
$ws.Range("E16").Value = "'  +2.75%  "
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.105.53"
$ws.Range("E17").Value = "'  +1.01%  "
$ws.Range("B18").Value = "'Polkadot"
$ws.Range("C18").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.21"
$ws.Range("E18").Value = "'  +2.44%  "
$ws.Range("D19").Value = "'13.02"
$ws.Range("E19").Value = "'  +0.01%  "
$ws.Range("D20").Value = "'8.21"
$ws.Range("E20").Value = "'  +0.61%  "
$ws.Range("D21").Value = "'344.04"
$ws.Range("E21").Value = "'  +1.59%  "
$ws.Range("E22").Value = "'  -0.19%  "
$ws.Range("E23").Value = "'  +1.91%  "
$ws.Range("D24").Value = "'65.84"
$ws.Range("E24").Value = "'  +0.39%  "
$ws.Range("E25").Value = "'  -0.62%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("D27").Value = "'0.0₃0937"
$ws.Range("E27").Value = "'  -0.58%  "
$ws.Range("E28").Value = "'  +4.39%  "
$ws.Range("D29").Value = "'7.30"
$ws.Range("E29").Value = "'  +2.59%  "
$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "'  +1.99%  "
$ws.Range("E31").Value = "'  +3.59%  "
$ws.Range("D32").Value = "'21.08"
$ws.Range("E32").Value = "'  +1.41%  "
$ws.Range("D33").Value = "'155.12"
$ws.Range("E33").Value = "'  +0.58%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("E34").Value = "'  +2.54%  "
$ws.Range("D35").Value = "'6.20"
$ws.Range("E35").Value = "'  +5.24%  "
$ws.Range("D36").Value = "'26.92"
$ws.Range("E36").Value = "'  +2.60%  "
$ws.Range("E37").Value = "'  +5.64%  "
$ws.Range("D38").Value = "'0.0688"
$ws.Range("E38").Value = "'  +1.27%  "
$ws.Range("D39").Value = "'3.96"
$ws.Range("E39").Value = "'  +2.44%  "
$ws.Range("D40").Value = "'3.145.43"
$ws.Range("E40").Value = "'  +1.05%  "
$ws.Range("D41").Value = "'36.88"
$ws.Range("E41").Value = "'  -0.43%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "'  -0.05%  "
$ws.Range("D43").Value = "'0.665"
$ws.Range("E43").Value = "'  -0.85%  "
$ws.Range("D44").Value = "'1.46"
$ws.Range("E44").Value = "'  +5.75%  "
$ws.Range("D45").Value = "'2.290.22"
$ws.Range("E45").Value = "'  +1.20%  "
$ws.Range("E46").Value = "'  +2.82%  "
$ws.Range("D47").Value = "'20.98"
$ws.Range("E47").Value = "'  +4.24%  "
$ws.Range("D48").Value = "'0.968"
$ws.Range("E48").Value = "'  +1.21%  "
$ws.Range("D49").Value = "'6.03"
$ws.Range("E49").Value = "'  +2.82%  "
$ws.Range("E50").Value = "'  +11.04%  "
$ws.Range("D51").Value = "'263.04"
$ws.Range("E51").Value = "'  +11.35%  "
